$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1. Locate the last paragraph of the document (currently the empty
#    list paragraph that ends the "Initialiseren van de sliders" list).
# ------------------------------------------------------------------
$lastIndex = $d.Paragraphs.Count
$lastPara  = $d.Paragraphs.Item($lastIndex)
$insertRange = $lastPara.Range
$insertRange.Collapse(0)   # wdCollapseEnd

# ------------------------------------------------------------------
# 2. Insert the new block of paragraphs (everything except the very
#    last, now-trailing, empty paragraph) right before the old last
#    paragraph's mark. Range.InsertXML always lands its content ahead
#    of the target paragraph mark, which keeps $lastPara intact (and
#    last) after the call - we fix $lastPara up afterwards.
# ------------------------------------------------------------------
$newBlockXml = @"
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">
  <pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">
    <pkg:xmlData>
      <w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">
        <w:body>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Lijstalinea"/>
              <w:numPr>
                <w:ilvl w:val="0"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Output wordt apart geregeld door </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>inline</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t xml:space="preserve"> </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>JavaScript</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
          </w:p>
          <w:p>
            <w:pPr>
              <w:pStyle w:val="Lijstalinea"/>
              <w:numPr>
                <w:ilvl w:val="1"/>
                <w:numId w:val="1"/>
              </w:numPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r><w:t xml:space="preserve">Value output moet ook wijzigen bij wijziging </w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Widht</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t>/</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Height</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t>&#8230;</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r><w:t>Zorgen voor start &#8216;</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>Player</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t>&#8217; met methode vanuit Class Board.</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
            <w:r><w:t>Ophalen van gegevens vanuit een dit m.b.v. data-</w:t></w:r>
            <w:proofErr w:type="spellStart"/>
            <w:r><w:t>value</w:t></w:r>
            <w:proofErr w:type="spellEnd"/>
            <w:r><w:t>=1;</w:t></w:r>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
          <w:p>
            <w:pPr>
              <w:spacing w:after="0"/>
            </w:pPr>
          </w:p>
        </w:body>
      </w:document>
    </pkg:xmlData>
  </pkg:part>
</pkg:package>
"@
$insertRange.InsertXML($newBlockXml)

# ------------------------------------------------------------------
# 3. The original last paragraph got pushed to the end again; strip
#    its list formatting (pStyle + numPr) so it becomes the plain
#    trailing empty paragraph the diff expects, keeping spacing after=0.
# ------------------------------------------------------------------
$trailingPara = $d.Paragraphs.Item($d.Paragraphs.Count)
$trailingPara.Style = "Standaard"
$trailingPara.Range.ParagraphFormat.SpaceAfter = 0

# ------------------------------------------------------------------
# 4. Turn the empty placeholder paragraph (the one right before the
#    trailing empty paragraph) into the hyperlink paragraph.
# ------------------------------------------------------------------
$hyperlinkParaIndex = $d.Paragraphs.Count - 1
$hyperlinkPara = $d.Paragraphs.Item($hyperlinkParaIndex)
$hyperlinkRange = $hyperlinkPara.Range
$hyperlinkRange.Collapse(0)
$d.Hyperlinks.Add($hyperlinkRange, "https://stackoverflow.com/questions/11238508/how-to-get-value-of-a-div-using-javascript") | Out-Null

# ------------------------------------------------------------------
# 5. Materialize + tune the "Hyperlink" character style so it matches
#    what a Dutch Word install records when a hyperlink is inserted.
# ------------------------------------------------------------------
$hyperlinkStyle = $d.Styles.Item("Hyperlink")
$hyperlinkStyle.BaseStyle = "Standaardalinea-lettertype"
$hyperlinkStyle.Priority = 99
$hyperlinkStyle.UnhideWhenUsed = $true
$hyperlinkStyle.QuickStyle = $false
$hyperlinkStyle.Font.TextColor.ObjectThemeColor = 10   # wdThemeColorHyperlink

# ------------------------------------------------------------------
# 6. Add the companion "Unresolved Mention" character style that Word
#    also records alongside Hyperlink in this document revision.
# ------------------------------------------------------------------
$mentionStyle = $d.Styles.Add("Onopgelostemelding", 2)   # wdStyleTypeCharacter
$mentionStyle.NameLocal = "Unresolved Mention"
$mentionStyle.BaseStyle = "Standaardalinea-lettertype"
$mentionStyle.Priority = 99
$mentionStyle.UnhideWhenUsed = $true
$mentionStyle.Font.TextColor.RGB = 6053472   # 0x605E5C (BGR-encoded)

Write-Host "Edit applied."
